# sample_selection_type_error.xlsx - add a new test product row (PROD7TEST)
# mirroring the existing "PROD6TEST" row (row 7), but with purchase_ok
# ("True") filled in and invoice_policy left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 7 (formatting + values) down into row 8 - two passes
# (formats, then values) so the new row keeps the same number formats /
# text-typed cells as the source row instead of Excel "smart" retyping
# strings like "true"/"11.11" as booleans/numbers.
$ws.Range("A7:H7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A7:H7").Copy()
$ws.Range("A8").PasteSpecial(-4163)

# purchase_ok was blank on row 7 - pull a real "True" text cell (F2) in so
# F8 ends up text-typed rather than being re-interpreted as a boolean.
$ws.Range("F2").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F2").Copy()
$ws.Range("F8").PasteSpecial(-4163)

# New product identity/name.
$ws.Range("A8").Value = "PROD7TEST"
$ws.Range("B8").Value = "Product 7 test"

# Leave the rest of row 8 (categ_id/type/sale_ok/standard_price/list_price)
# identical to row 7 - already copied above.

$ws.Range("F17").Select()
